# Update overall stats on the active worksheet.
# Rows 2-6 all share the same updated values for columns B, C, E, G, I, J, K, M, N, O, P.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values keyed by column letter
$values = @{
    "B" = 1
    "C" = 1
    "E" = 241.5
    "G" = 15.181125
    "I" = 5.25
    "J" = 226.32
    "K" = 226.32
    "M" = 0.5
    "N" = 0.5
    "O" = 226.32
    "P" = 15.91
}

foreach ($row in 2..6) {
    foreach ($col in $values.Keys) {
        $ws.Range("$col$row").Value = $values[$col]
    }
}
